# Auto-generated Excel COM-interop script to update cryptos.xlsx price/volume data
# (mirrors the GitHub Actions crypto-data refresh commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.946.63'
$ws.Range("E2").Value = '  -0.13%  '
$ws.Range("D3").Value = '1.673.49'
$ws.Range("E3").Value = '  +1.05%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("E6").Value = '  +1.50%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("E9").Value = '  +0.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.19'
$ws.Range("E10").Value = '  +0.17%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0889'
$ws.Range("E11").Value = '  +1.17%  '
$ws.Range("D12").Value = '1.908.63'
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("D13").Value = '1.666.26'
$ws.Range("E13").Value = '  +0.65%  '
$ws.Range("E14").Value = '  -0.04%  '
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.68'
$ws.Range("D17").Value = '26.945.86'
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("E18").Value = '  +3.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '234.43'
$ws.Range("E19").Value = '  -1.09%  '
$ws.Range("D20").Value = '0.0₃0733'
$ws.Range("E20").Value = '  -1.09%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("E23").Value = '  -1.39%  '
$ws.Range("E24").Value = '  -1.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.56'
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.15'
$ws.Range("E26").Value = '  +0.30%  '
$ws.Range("E27").Value = '  +0.63%  '
$ws.Range("E28").Value = '  -1.39%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  -0.16%  '
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("E32").Value = '  +0.81%  '
$ws.Range("D33").Value = '1.476.27'
$ws.Range("E33").Value = '  -5.03%  '
$ws.Range("E34").Value = '  +2.12%  '
$ws.Range("E35").Value = '  +2.64%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.42'
$ws.Range("E36").Value = '  +0.20%  '
$ws.Range("E37").Value = '  -0.51%  '
$ws.Range("E38").Value = '  -1.09%  '
$ws.Range("E39").Value = '  +0.83%  '
$ws.Range("E40").Value = '  +8.27%  '
$ws.Range("E41").Value = '  -3.73%  '
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("E43").Value = '  +2.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '66.72'
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").Value = '1.815.30'
$ws.Range("E45").Value = '  +1.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.779'
$ws.Range("E46").Value = '  +0.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.42'
$ws.Range("E47").Value = '  +0.20%  '
$ws.Range("E48").Value = '  -0.16%  '
$ws.Range("E49").Value = '  +1.58%  '
$ws.Range("E50").Value = '  +0.44%  '
$ws.Range("E51").Value = '  +0.62%  '
